$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Rushing": C.Wentz no longer has stats this week - remove
# his row (shifting everyone up), then append J.Doyle's Week 16
# rushing line. A few other players' cumulative stats also bumped.
# ---------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

$rushing.Rows("2").Delete()

# Updated cumulative totals after Week 16
$rushing.Range("C3").Value = 160
$rushing.Range("D3").Value = 105
$rushing.Range("E3").Value = 31
$rushing.Range("F3").Value = 79

$rushing.Range("C4").Value = 26
$rushing.Range("D4").Value = 21

$rushing.Range("D7").Value = 2

# New row for J.Doyle
$rushing.Range("A9").Copy($rushing.Range("A10"))
$rushing.Range("A10").Value = 9
$rushing.Range("B10").Value = "J.Doyle"
$rushing.Range("C10").Value = 1
$rushing.Range("D10").Value = 0
$rushing.Range("E10").Value = 0
$rushing.Range("F10").Value = 1
$rushing.Range("D10").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "Receiving": update cumulative totals after Week 16 (row
# order unchanged).
# ---------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

$receiving.Range("C2").Value = 43

$receiving.Range("C3").Value = 44
$receiving.Range("D3").Value = 31

$receiving.Range("C5").Value = 93
$receiving.Range("D5").Value = 66
$receiving.Range("E5").Value = 23
$receiving.Range("F5").Value = 11
$receiving.Range("G5").Value = 16
$receiving.Range("H5").Value = 7

$receiving.Range("C9").Value = 15
$receiving.Range("E9").Value = 5
$receiving.Range("F9").Value = 2

$receiving.Range("C10").Value = 36
$receiving.Range("D10").Value = 30
$receiving.Range("E10").Value = 7
$receiving.Range("F10").Value = 4
$receiving.Range("G10").Value = 3
$receiving.Range("H10").Value = 3

$receiving.Range("C11").Value = 3
$receiving.Range("D11").Value = 2
$receiving.Range("G11").Value = 1
$receiving.Range("H11").Value = 1

$receiving.Range("C14").Value = 27
$receiving.Range("D14").Value = 16
$receiving.Range("E14").Value = 13
$receiving.Range("F14").Value = 5
$receiving.Range("G14").Value = 9
$receiving.Range("H14").Value = 5

$receiving.Range("C15").Value = 13
$receiving.Range("D15").Value = 10
